$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.73999977111816
$ws.Range("C2").Value = 19.84000015258789
$ws.Range("D2").Value = 22.50360648358454
